$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("48:48").Insert()

$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 44533
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100109
$ws.Range("H48").Value = "Uva"
$ws.Range("I48").Value = 100109001
$ws.Range("J48").Value = "Uva"
$ws.Range("K48").Value = "Superior Seedless"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 300
$ws.Range("N48").Value = 27000
$ws.Range("O48").Value = 28000
$ws.Range("P48").Value = 27500
$ws.Range("Q48").Value = "$/bandeja 8 kilos"
$ws.Range("R48").Value = "Provincia de Limarí"
$ws.Range("S48").Value = 3438
$ws.Range("T48").Value = 8

$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
